$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row (row 33) by duplicating the formatting of the last data
# row (32), then overwriting the values with the new user's data. This
# keeps per-column cell styles (e.g. the left-aligned email column)
# consistent with the rest of the table.
$ws.Range("A32:K32").Copy($ws.Range("A33:K33"))

$ws.Range("A33").Value = 110032
$ws.Range("B33").Value = 9317596770
$ws.Range("C33").Value = "Ewan Marsh"
$ws.Range("D33").Value = "ewan.marsh@xyz.com"
$ws.Range("E33").Value = 818876433
$ws.Range("F33").Value = "ACT"
$ws.Range("G33").Value = "eng"
$ws.Range("H33").Value = "PWD"
$ws.Range("I33").Value = $true
$ws.Range("J33").Value = "superadmin"
$ws.Range("K33").Value = "now()"

# Select column L (whole column, to the last column) with L1 as the
# active cell, matching the saved view state.
$ws.Range("L1:XFD1048576").Select()

# Page setup: vertical print quality now specified (300 dpi) instead of 0.
$ws.PageSetup.PrintQuality = 300

$wb.Save()
